$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F3").Value = 1190
$ws.Range("F4").Value = 2658
$ws.Range("F5").Value = 1190
$ws.Range("F6").Value = 2658
